$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1928104575163399
$ws.Range("C2").Value = 0.5522875816993464
$ws.Range("J2").Value = 0.0261437908496732
$ws.Range("P2").Value = 0.1470588235294118
$ws.Range("S2").Value = 0.08169934640522876
$ws.Range("B3").Value = 0.01052631578947368
$ws.Range("C3").Value = 0.04736842105263158
$ws.Range("J3").Value = 0.01578947368421053
$ws.Range("P3").Value = 0.7947368421052632
$ws.Range("S3").Value = 0.131578947368421
$ws.Range("J4").Value = 0.09375
$ws.Range("P4").Value = 0.6875
$ws.Range("S4").Value = 0.21875
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.07421875
$ws.Range("D6").Value = 0.01953125
$ws.Range("E6").Value = 0.00390625
$ws.Range("F6").Value = 0.07421875
$ws.Range("J6").Value = 0.25
$ws.Range("O6").Value = 0.015625
$ws.Range("Q6").Value = 0.171875
$ws.Range("R6").Value = 0.05078125
$ws.Range("S6").Value = 0.33984375
$ws.Range("B7").Value = 0.09345794392523364
$ws.Range("D7").Value = 0.01869158878504673
$ws.Range("F7").Value = 0.06542056074766354
$ws.Range("J7").Value = 0.1308411214953271
$ws.Range("O7").Value = 0.01401869158878505
$ws.Range("Q7").Value = 0.2009345794392523
$ws.Range("R7").Value = 0.08411214953271028
$ws.Range("S7").Value = 0.3925233644859813
$ws.Range("B8").Value = 0.09411764705882353
$ws.Range("D8").Value = 0.02745098039215686
$ws.Range("E8").Value = 0.00196078431372549
$ws.Range("F8").Value = 0.06470588235294118
$ws.Range("J8").Value = 0.1274509803921569
$ws.Range("O8").Value = 0.02549019607843137
$ws.Range("Q8").Value = 0.2
$ws.Range("R8").Value = 0.1215686274509804
$ws.Range("S8").Value = 0.3372549019607843
$ws.Range("B9").Value = 0.0625
$ws.Range("D9").Value = 0.02083333333333333
$ws.Range("F9").Value = 0.04166666666666666
$ws.Range("J9").Value = 0.1354166666666667
$ws.Range("O9").Value = 0.02083333333333333
$ws.Range("Q9").Value = 0.1354166666666667
$ws.Range("R9").Value = 0.1354166666666667
$ws.Range("S9").Value = 0.4479166666666667
$ws.Range("B10").Value = 0.1020999275887038
$ws.Range("D10").Value = 0.02751629254163649
$ws.Range("E10").Value = 0.001448225923244026
$ws.Range("F10").Value = 0.06879073135409124
$ws.Range("J10").Value = 0.1245474293989862
$ws.Range("O10").Value = 0.01448225923244026
$ws.Range("Q10").Value = 0.2498189717595945
$ws.Range("R10").Value = 0.09558291093410572
$ws.Range("S10").Value = 0.3157132512671977
$ws.Range("F11").Value = 0.006097560975609756
$ws.Range("G11").Value = 0.1432926829268293
$ws.Range("J11").Value = 0.0975609756097561
$ws.Range("K11").Value = 0.1981707317073171
$ws.Range("L11").Value = 0.5335365853658537
$ws.Range("S11").Value = 0.02134146341463415
$ws.Range("G12").Value = 0.8057142857142857
$ws.Range("J12").Value = 0.1542857142857143
$ws.Range("K12").Value = 0.005714285714285714
$ws.Range("L12").Value = 0.01714285714285714
$ws.Range("S12").Value = 0.01714285714285714
$ws.Range("G13").Value = 0.6888888888888889
$ws.Range("J13").Value = 0.2444444444444444
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("F15").Value = 0.03238866396761134
$ws.Range("H15").Value = 0.1376518218623482
$ws.Range("I15").Value = 0.03238866396761134
$ws.Range("J15").Value = 0.3967611336032389
$ws.Range("K15").Value = 0.06477732793522267
$ws.Range("M15").Value = 0.0242914979757085
$ws.Range("N15").Value = 0.004048582995951417
$ws.Range("O15").Value = 0.04453441295546558
$ws.Range("S15").Value = 0.2631578947368421
$ws.Range("F16").Value = 0.03097345132743363
$ws.Range("H16").Value = 0.1504424778761062
$ws.Range("I16").Value = 0.06637168141592921
$ws.Range("J16").Value = 0.411504424778761
$ws.Range("K16").Value = 0.1194690265486726
$ws.Range("M16").Value = 0.008849557522123894
$ws.Range("O16").Value = 0.06637168141592921
$ws.Range("S16").Value = 0.1460176991150443
$ws.Range("F17").Value = 0.02692998204667863
$ws.Range("H17").Value = 0.177737881508079
$ws.Range("I17").Value = 0.0466786355475763
$ws.Range("J17").Value = 0.4344703770197487
$ws.Range("K17").Value = 0.1166965888689408
$ws.Range("M17").Value = 0.01436265709156194
$ws.Range("O17").Value = 0.07360861759425494
$ws.Range("S17").Value = 0.1095152603231598
$ws.Range("F18").Value = 0.02092050209205021
$ws.Range("H18").Value = 0.2092050209205021
$ws.Range("I18").Value = 0.01255230125523013
$ws.Range("J18").Value = 0.4184100418410042
$ws.Range("K18").Value = 0.100418410041841
$ws.Range("M18").Value = 0.01673640167364017
$ws.Range("N18").Value = 0.004184100418410041
$ws.Range("O18").Value = 0.07949790794979079
$ws.Range("S18").Value = 0.1380753138075314
$ws.Range("F19").Value = 0.02333333333333333
$ws.Range("H19").Value = 0.2425
$ws.Range("I19").Value = 0.03416666666666666
$ws.Range("J19").Value = 0.3758333333333334
$ws.Range("K19").Value = 0.1075
$ws.Range("M19").Value = 0.02083333333333333
$ws.Range("O19").Value = 0.07333333333333333
$ws.Range("S19").Value = 0.1225
